# Auto-generated Excel COM-interop script to apply the diff to cryptos.xlsx
# Updates cell text values for the crypto price table (rows 2-51, columns B-E)
# All target cells are text strings; we force text storage and then clear the
# format override so no stray cell style is introduced (matches original styling).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextCell $ws.Range("D2") '25.639.19'
Set-TextCell $ws.Range("E2") '  -3.25%  '
Set-TextCell $ws.Range("D3") '1.738.83'
Set-TextCell $ws.Range("E3") '  -5.60%  '
Set-TextCell $ws.Range("E4") '  +0.09%  '
Set-TextCell $ws.Range("D5") '235.51'
Set-TextCell $ws.Range("E5") '  -10.33%  '
Set-TextCell $ws.Range("E6") '  -0.01%  '
Set-TextCell $ws.Range("D7") '0.4892'
Set-TextCell $ws.Range("E7") '  -8.07%  '
Set-TextCell $ws.Range("D8") '41.30'
Set-TextCell $ws.Range("E8") '  -8.35%  '
Set-TextCell $ws.Range("D9") '0.2547'
Set-TextCell $ws.Range("E9") '  -17.54%  '
Set-TextCell $ws.Range("D10") '0.06069'
Set-TextCell $ws.Range("E10") '  -12.01%  '
Set-TextCell $ws.Range("D11") '1.741.55'
Set-TextCell $ws.Range("E11") '  -5.47%  '
Set-TextCell $ws.Range("D12") '0.06834'
Set-TextCell $ws.Range("E12") '  -12.73%  '
Set-TextCell $ws.Range("D13") '14.74'
Set-TextCell $ws.Range("E13") '  -20.56%  '
Set-TextCell $ws.Range("D14") '4.433'
Set-TextCell $ws.Range("E14") '  -12.14%  '
Set-TextCell $ws.Range("D15") '75.62'
Set-TextCell $ws.Range("E15") '  -15.67%  '
Set-TextCell $ws.Range("D16") '0.5581'
Set-TextCell $ws.Range("E16") '  -26.65%  '
Set-TextCell $ws.Range("E17") '  -0.04%  '
Set-TextCell $ws.Range("E18") '  -0.02%  '
Set-TextCell $ws.Range("D19") '25.677.28'
Set-TextCell $ws.Range("E19") '  -3.23%  '
Set-TextCell $ws.Range("D20") '11.44'
Set-TextCell $ws.Range("E20") '  -18.46%  '
Set-TextCell $ws.Range("D21") '0.000006532'
Set-TextCell $ws.Range("E21") '  -17.83%  '
Set-TextCell $ws.Range("D22") '1.961.51'
Set-TextCell $ws.Range("E22") '  -5.53%  '
Set-TextCell $ws.Range("D23") '4.025'
Set-TextCell $ws.Range("E23") '  -13.08%  '
Set-TextCell $ws.Range("D24") '7.887'
Set-TextCell $ws.Range("E24") '  -15.37%  '
Set-TextCell $ws.Range("D25") '4.986'
Set-TextCell $ws.Range("E25") '  -17.04%  '
Set-TextCell $ws.Range("D26") '137.13'
Set-TextCell $ws.Range("E26") '  -3.25%  '
Set-TextCell $ws.Range("D27") '1.475'
Set-TextCell $ws.Range("E27") '  -12.67%  '
Set-TextCell $ws.Range("D28") '1.819'
Set-TextCell $ws.Range("E28") '  -16.97%  '
Set-TextCell $ws.Range("E29") '  -13.84%  '
Set-TextCell $ws.Range("D30") '100.88'
Set-TextCell $ws.Range("E30") '  -9.22%  '
Set-TextCell $ws.Range("E31") '  -9.61%  '
Set-TextCell $ws.Range("D32") '3.674'
Set-TextCell $ws.Range("E32") '  -14.19%  '
Set-TextCell $ws.Range("D33") '3.363'
Set-TextCell $ws.Range("E33") '  -17.90%  '
Set-TextCell $ws.Range("D34") '0.04399'
Set-TextCell $ws.Range("E34") '  -8.90%  '
Set-TextCell $ws.Range("D35") '1.000'
Set-TextCell $ws.Range("E35") '  +0.02%  '
Set-TextCell $ws.Range("D36") '2.610'
Set-TextCell $ws.Range("E36") '  -11.05%  '
Set-TextCell $ws.Range("D37") '0.9674'
Set-TextCell $ws.Range("E37") '  -14.71%  '
Set-TextCell $ws.Range("D38") '0.5868'
Set-TextCell $ws.Range("E38") '  -19.97%  '
Set-TextCell $ws.Range("D39") '2.647'
Set-TextCell $ws.Range("E39") '  -14.71%  '
Set-TextCell $ws.Range("D41") '103.19'
Set-TextCell $ws.Range("E41") '  -4.65%  '
Set-TextCell $ws.Range("B42") 'PaxosStandard'
Set-TextCell $ws.Range("C42") 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
Set-TextCell $ws.Range("D42") '1.002'
Set-TextCell $ws.Range("E42") '  +0.03%  '
Set-TextCell $ws.Range("B43") 'VeChain'
Set-TextCell $ws.Range("C43") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range("D43") '0.01497'
Set-TextCell $ws.Range("E43") '  -13.18%  '
Set-TextCell $ws.Range("B44") 'RenderToken'
Set-TextCell $ws.Range("C44") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range("D44") '1.864'
Set-TextCell $ws.Range("E44") '  -19.70%  '
Set-TextCell $ws.Range("B45") 'FraxShare'
Set-TextCell $ws.Range("C45") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range("D45") '5.124'
Set-TextCell $ws.Range("E45") '  -13.07%  '
Set-TextCell $ws.Range("B46") 'TheSandbox'
Set-TextCell $ws.Range("C46") 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws.Range("D46") '0.3712'
Set-TextCell $ws.Range("E46") '  -22.79%  '
Set-TextCell $ws.Range("B47") 'TrustWalletToken'
Set-TextCell $ws.Range("C47") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws.Range("D47") '0.7230'
Set-TextCell $ws.Range("E47") '  -19.91%  '
Set-TextCell $ws.Range("B48") 'Cronos'
Set-TextCell $ws.Range("C48") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range("D48") '0.05210'
Set-TextCell $ws.Range("E48") '  -10.26%  '
Set-TextCell $ws.Range("B49") 'Algorand'
Set-TextCell $ws.Range("C49") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws.Range("D49") '0.1076'
Set-TextCell $ws.Range("E49") '  -13.29%  '
Set-TextCell $ws.Range("B50") 'Elrond'
Set-TextCell $ws.Range("C50") 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell $ws.Range("D50") '29.80'
Set-TextCell $ws.Range("E50") '  -14.81%  '
Set-TextCell $ws.Range("D51") '51.59'
Set-TextCell $ws.Range("E51") '  -14.55%  '
